$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Tanana"
$ws.Cells.Item(2, 3).Value = "Mariam"

$ws.Cells.Item(3, 2).Value = "El Haddad"
$ws.Cells.Item(3, 3).Value = "Mohamed"

$ws.Cells.Item(4, 3).Value = "Khalid"
$ws.Cells.Item(4, 2).Value = "Amechnoue"

$ws.Cells.Item(5, 2).Value = "Massou"
$ws.Cells.Item(5, 3).Value = "Siham"

$ws.Cells.Item(6, 2).Value = "Belmokadem"

$ws.Cells.Item(7, 2).Value = "Ghailani"
$ws.Cells.Item(7, 3).Value = "Mohamed"

$ws.Cells.Item(8, 2).Value = "Fissoune"
$ws.Cells.Item(8, 3).Value = "Rachida"

$ws.Cells.Item(9, 2).Value = "Nait Bouker"
$ws.Cells.Item(9, 3).Value = "Nezha"

$ws.Range("C9").Select()
